# Qubit on CFX: when both HS and BR assays are measured for a sample
# (column F = "both"), the script chooses which measurement (HS "C" or
# BR "E") should be used for the merged [DNA] ng/uL column (G).
# Fill in the previously-empty results for rows 21, 25, 30, 36 and 37.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G21").Value = 2.464824779166207
$ws.Range("G25").Value = 4.170888231044157
$ws.Range("G30").Value = 4.302692002788298
$ws.Range("G36").Value = 3.871860744388574
$ws.Range("G37").Value = 5.11186148980102
